$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Shared-string text edit: "Field:id" -> "Field:HAS" (cell B1, reserved key)
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Field:HAS"

# ---------------------------------------------------------------------------
# 2. New column E (dimension grows from A1:D10 to A1:E10). Column E gets the
#    same width as the rest of the unbounded range (col F onward) it used to
#    share, but now carries its own bordered "box" look:
#      E1      -> left/right/top thin red border   (top edge of the box)
#      E2:E9   -> left/right thin red border only   (sides of the box)
#      E10     -> left/right/bottom thin red border (bottom edge of the box)
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = $ws.Columns("F").ColumnWidth

# --- E1: copy D1's fill/format, then carve out the top-of-box border -------
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Borders.Item(7).LineStyle = $ws.Range("D1").Borders.Item(7).LineStyle
$ws.Range("E1").Borders.Item(7).Color = $ws.Range("D1").Borders.Item(7).Color
$ws.Range("E1").Borders.Item(7).Weight = $ws.Range("D1").Borders.Item(7).Weight
$ws.Range("E1").Borders.Item(10).LineStyle = $ws.Range("D1").Borders.Item(10).LineStyle
$ws.Range("E1").Borders.Item(10).Color = $ws.Range("D1").Borders.Item(10).Color
$ws.Range("E1").Borders.Item(10).Weight = $ws.Range("D1").Borders.Item(10).Weight
$ws.Range("E1").Borders.Item(8).LineStyle = $ws.Range("D1").Borders.Item(8).LineStyle
$ws.Range("E1").Borders.Item(8).Color = $ws.Range("D1").Borders.Item(8).Color
$ws.Range("E1").Borders.Item(8).Weight = $ws.Range("D1").Borders.Item(8).Weight
$ws.Range("E1").Borders.Item(9).LineStyle = -4142
$ws.Range("E1").VerticalAlignment = -4107

# --- E2: copy D2's fill/format, then carve out the side-of-box border ------
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E2").Borders.Item(7).LineStyle = $ws.Range("E1").Borders.Item(7).LineStyle
$ws.Range("E2").Borders.Item(7).Color = $ws.Range("E1").Borders.Item(7).Color
$ws.Range("E2").Borders.Item(7).Weight = $ws.Range("E1").Borders.Item(7).Weight
$ws.Range("E2").Borders.Item(10).LineStyle = $ws.Range("E1").Borders.Item(10).LineStyle
$ws.Range("E2").Borders.Item(10).Color = $ws.Range("E1").Borders.Item(10).Color
$ws.Range("E2").Borders.Item(10).Weight = $ws.Range("E1").Borders.Item(10).Weight
$ws.Range("E2").Borders.Item(8).LineStyle = -4142
$ws.Range("E2").Borders.Item(9).LineStyle = -4142
$ws.Range("E2").VerticalAlignment = -4107

# --- Propagate E2's "sides only" look down through E9 in one shot ----------
$ws.Range("E2").Copy()
$ws.Range("E3:E9").PasteSpecial(-4122)

# --- E10: copy D10's fill/format, then carve out the bottom-of-box border --
$ws.Range("D10").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Borders.Item(7).LineStyle = $ws.Range("E2").Borders.Item(7).LineStyle
$ws.Range("E10").Borders.Item(7).Color = $ws.Range("E2").Borders.Item(7).Color
$ws.Range("E10").Borders.Item(7).Weight = $ws.Range("E2").Borders.Item(7).Weight
$ws.Range("E10").Borders.Item(10).LineStyle = $ws.Range("E2").Borders.Item(10).LineStyle
$ws.Range("E10").Borders.Item(10).Color = $ws.Range("E2").Borders.Item(10).Color
$ws.Range("E10").Borders.Item(10).Weight = $ws.Range("E2").Borders.Item(10).Weight
$ws.Range("E10").Borders.Item(8).LineStyle = -4142
$ws.Range("E10").Borders.Item(9).LineStyle = $ws.Range("E1").Borders.Item(7).LineStyle
$ws.Range("E10").Borders.Item(9).Color = $ws.Range("E1").Borders.Item(7).Color
$ws.Range("E10").Borders.Item(9).Weight = $ws.Range("E1").Borders.Item(7).Weight
$ws.Range("E10").VerticalAlignment = -4107

Write-Output "edit applied"
